$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header of column A from "Gen" to "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 2. Replace the "Gen" values in column A (generation counts) with the
#    fractional "MaxFES" values used by the new PSO-based results.
$maxfesValues = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $maxfesValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $maxfesValues[$i]
}

# 3. Drop the "Run 50" column (column AZ), shrinking the table from 51 runs
#    (Run 0..Run 50) down to 50 runs (Run 0..Run 49). The old "Mean" column
#    (BA) then shifts left to become the new last column (AZ).
$ws.Columns("AZ").Delete()

# 4. Recompute the Mean column (now column AZ, after the deletion above)
#    for each remaining row, using the new set of 50 run columns (B..AY).
$meanValues = @(13.71434814, 12.94602389, 10.36595172, 6.89565865, 5.7401227, 5.20886074, 4.79054858, 4.39546716, 4.08847069, 3.78573622, 3.57405791, 3.33237497, 3.17170225)
for ($i = 0; $i -lt $meanValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 52).Value = $meanValues[$i]
}
